# Apply a re-ordering of several data rows in the "Artfynd" sheet.
# The underlying report was re-exported with a different row ordering for a
# handful of records (rows 11-13, 16-17 and 19-21); the cell content for
# every column (A:AY) simply moved to a different row number. We snapshot
# the full row ranges first (so the subsequent writes can't clobber a value
# we still need to read) and then write each snapshot to its new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol  = "AY"

function Get-RowValues($sheet, $rowNum) {
    return $sheet.Range("$firstCol$rowNum`:$lastCol$rowNum").Value2
}

# Snapshot every source row BEFORE any writes happen.
$row11 = Get-RowValues $ws 11
$row12 = Get-RowValues $ws 12
$row13 = Get-RowValues $ws 13
$row16 = Get-RowValues $ws 16
$row17 = Get-RowValues $ws 17
$row19 = Get-RowValues $ws 19
$row20 = Get-RowValues $ws 20
$row21 = Get-RowValues $ws 21

# Destination row <- source row content
# 11 <- 13, 12 <- 11, 13 <- 12
# 16 <- 17, 17 <- 16
# 19 <- 21, 20 <- 19, 21 <- 20
$ws.Range("$firstCol" + "11:" + "$lastCol" + "11").Value2 = $row13
$ws.Range("$firstCol" + "12:" + "$lastCol" + "12").Value2 = $row11
$ws.Range("$firstCol" + "13:" + "$lastCol" + "13").Value2 = $row12

$ws.Range("$firstCol" + "16:" + "$lastCol" + "16").Value2 = $row17
$ws.Range("$firstCol" + "17:" + "$lastCol" + "17").Value2 = $row16

$ws.Range("$firstCol" + "19:" + "$lastCol" + "19").Value2 = $row21
$ws.Range("$firstCol" + "20:" + "$lastCol" + "20").Value2 = $row19
$ws.Range("$firstCol" + "21:" + "$lastCol" + "21").Value2 = $row20
